$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics (ligand/receptor expression & specificity scores)
# for the Pgf-Flt1 ligand-receptor pair sheet, per updated script run.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.961753666666667
$ws.Range("H2").Value = 17.885261
$ws.Range("I2").Value = 0.7691652453336842
$ws.Range("J2").Value = 0.7691652453336842
$ws.Range("M2").Value = 154.942487
$ws.Range("N2").Value = 464.827461
$ws.Range("O2").Value = 0.982851703624775
$ws.Range("P2").Value = 0.9828517036247751
$ws.Range("Q2").Value = 923.7289399947023
$ws.Range("R2").Value = 8313.560459952321
$ws.Range("S2").Value = 0.7559753717451796
$ws.Range("T2").Value = 0.7559753717451797

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.961753666666667
$ws.Range("H3").Value = 17.885261
$ws.Range("I3").Value = 0.7691652453336842
$ws.Range("J3").Value = 0.7691652453336842
$ws.Range("O3").Value = 0.003358739549735124
$ws.Range("P3").Value = 0.003358739549735124
$ws.Range("Q3").Value = 3.156696897968222
$ws.Range("R3").Value = 28.410272081714
$ws.Range("S3").Value = 0.002583425729783965
$ws.Range("T3").Value = 0.002583425729783965

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.961753666666667
$ws.Range("H4").Value = 17.885261
$ws.Range("I4").Value = 0.7691652453336842
$ws.Range("J4").Value = 0.7691652453336842
$ws.Range("M4").Value = 1.771368666666667
$ws.Range("N4").Value = 5.314106000000001
$ws.Range("O4").Value = 0.01123638032078883
$ws.Range("P4").Value = 0.01123638032078884
$ws.Range("Q4").Value = 10.56046364351845
$ws.Range("R4").Value = 95.04417279166601
$ws.Range("S4").Value = 0.008642633226102125
$ws.Range("T4").Value = 0.008642633226102127

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.961753666666667
$ws.Range("H5").Value = 17.885261
$ws.Range("I5").Value = 0.7691652453336842
$ws.Range("J5").Value = 0.7691652453336842
$ws.Range("M5").Value = 0.4024976666666666
$ws.Range("N5").Value = 1.207493
$ws.Range("O5").Value = 0.002553176504700935
$ws.Range("P5").Value = 0.002553176504700936
$ws.Range("Q5").Value = 2.399591940074778
$ws.Range("R5").Value = 21.596327460673
$ws.Range("S5").Value = 0.001963814632618494
$ws.Range("T5").Value = 0.001963814632618494

# Row 6
$ws.Range("I6").Value = 0.07350643921898506
$ws.Range("J6").Value = 0.07350643921898506
$ws.Range("M6").Value = 154.942487
$ws.Range("N6").Value = 464.827461
$ws.Range("O6").Value = 0.982851703624775
$ws.Range("P6").Value = 0.9828517036247751
$ws.Range("Q6").Value = 88.277552313328
$ws.Range("R6").Value = 794.497970819952
$ws.Range("S6").Value = 0.07224592901377044
$ws.Range("T6").Value = 0.07224592901377046

# Row 7
$ws.Range("I7").Value = 0.07350643921898506
$ws.Range("J7").Value = 0.07350643921898506
$ws.Range("O7").Value = 0.003358739549735124
$ws.Range("P7").Value = 0.003358739549735124
$ws.Range("S7").Value = 0.0002468889845650061
$ws.Range("T7").Value = 0.0002468889845650061

# Row 8
$ws.Range("I8").Value = 0.07350643921898506
$ws.Range("J8").Value = 0.07350643921898506
$ws.Range("M8").Value = 1.771368666666667
$ws.Range("N8").Value = 5.314106000000001
$ws.Range("O8").Value = 0.01123638032078883
$ws.Range("P8").Value = 0.01123638032078884
$ws.Range("Q8").Value = 1.009226669621333
$ws.Range("R8").Value = 9.083040026592002
$ws.Range("S8").Value = 0.0008259463070914643
$ws.Range("T8").Value = 0.0008259463070914644

# Row 9
$ws.Range("I9").Value = 0.07350643921898506
$ws.Range("J9").Value = 0.07350643921898506
$ws.Range("M9").Value = 0.4024976666666666
$ws.Range("N9").Value = 1.207493
$ws.Range("O9").Value = 0.002553176504700935
$ws.Range("P9").Value = 0.002553176504700936
$ws.Range("Q9").Value = 0.2293206305973333
$ws.Range("R9").Value = 2.063885675376
$ws.Range("S9").Value = 0.00018767491355814
$ws.Range("T9").Value = 0.0001876749135581401

# Row 10
$ws.Range("G10").Value = 1.143196
$ws.Range("H10").Value = 3.429588
$ws.Range("I10").Value = 0.1474912720263607
$ws.Range("J10").Value = 0.1474912720263607
$ws.Range("M10").Value = 154.942487
$ws.Range("N10").Value = 464.827461
$ws.Range("O10").Value = 0.982851703624775
$ws.Range("P10").Value = 0.9828517036247751
$ws.Range("Q10").Value = 177.129631368452
$ws.Range("R10").Value = 1594.166682316068
$ws.Range("S10").Value = 0.1449620479808937
$ws.Range("T10").Value = 0.1449620479808937

# Row 11
$ws.Range("G11").Value = 1.143196
$ws.Range("H11").Value = 3.429588
$ws.Range("I11").Value = 0.1474912720263607
$ws.Range("J11").Value = 0.1474912720263607
$ws.Range("O11").Value = 0.003358739549735124
$ws.Range("P11").Value = 0.003358739549735124
$ws.Range("Q11").Value = 0.6053123743013332
$ws.Range("R11").Value = 5.447811368711999
$ws.Range("S11").Value = 0.0004953847685956794
$ws.Range("T11").Value = 0.0004953847685956794

# Row 12
$ws.Range("G12").Value = 1.143196
$ws.Range("H12").Value = 3.429588
$ws.Range("I12").Value = 0.1474912720263607
$ws.Range("J12").Value = 0.1474912720263607
$ws.Range("M12").Value = 1.771368666666667
$ws.Range("N12").Value = 5.314106000000001
$ws.Range("O12").Value = 0.01123638032078883
$ws.Range("P12").Value = 0.01123638032078884
$ws.Range("Q12").Value = 2.025021574258667
$ws.Range("R12").Value = 18.225194168328
$ws.Range("S12").Value = 0.001657268026485112
$ws.Range("T12").Value = 0.001657268026485112

# Row 13
$ws.Range("G13").Value = 1.143196
$ws.Range("H13").Value = 3.429588
$ws.Range("I13").Value = 0.1474912720263607
$ws.Range("J13").Value = 0.1474912720263607
$ws.Range("M13").Value = 0.4024976666666666
$ws.Range("N13").Value = 1.207493
$ws.Range("O13").Value = 0.002553176504700935
$ws.Range("P13").Value = 0.002553176504700936
$ws.Range("Q13").Value = 0.4601337225426666
$ws.Range("R13").Value = 4.141203502883999
$ws.Range("S13").Value = 0.0003765712503861584
$ws.Range("T13").Value = 0.0003765712503861584

# Row 14
$ws.Range("G14").Value = 0.07624633333333333
$ws.Range("H14").Value = 0.228739
$ws.Range("I14").Value = 0.009837043420970016
$ws.Range("J14").Value = 0.009837043420970016
$ws.Range("M14").Value = 154.942487
$ws.Range("N14").Value = 464.827461
$ws.Range("O14").Value = 0.982851703624775
$ws.Range("P14").Value = 0.9828517036247751
$ws.Range("Q14").Value = 11.81379651129767
$ws.Range("R14").Value = 106.324168601679
$ws.Range("S14").Value = 0.009668354884931266
$ws.Range("T14").Value = 0.009668354884931266

# Row 15
$ws.Range("G15").Value = 0.07624633333333333
$ws.Range("H15").Value = 0.228739
$ws.Range("I15").Value = 0.009837043420970016
$ws.Range("J15").Value = 0.009837043420970016
$ws.Range("O15").Value = 0.003358739549735124
$ws.Range("P15").Value = 0.003358739549735124
$ws.Range("Q15").Value = 0.04037177269844444
$ws.Range("R15").Value = 0.363345954286
$ws.Range("S15").Value = 0.00003304006679047369
$ws.Range("T15").Value = 0.00003304006679047369

# Row 16
$ws.Range("G16").Value = 0.07624633333333333
$ws.Range("H16").Value = 0.228739
$ws.Range("I16").Value = 0.009837043420970016
$ws.Range("J16").Value = 0.009837043420970016
$ws.Range("M16").Value = 1.771368666666667
$ws.Range("N16").Value = 5.314106000000001
$ws.Range("O16").Value = 0.01123638032078883
$ws.Range("P16").Value = 0.01123638032078884
$ws.Range("Q16").Value = 0.1350603658148889
$ws.Range("R16").Value = 1.215543292334
$ws.Range("S16").Value = 0.0001105327611101328
$ws.Range("T16").Value = 0.0001105327611101328

# Row 17
$ws.Range("G17").Value = 0.07624633333333333
$ws.Range("H17").Value = 0.228739
$ws.Range("I17").Value = 0.009837043420970016
$ws.Range("J17").Value = 0.009837043420970016
$ws.Range("M17").Value = 0.4024976666666666
$ws.Range("N17").Value = 1.207493
$ws.Range("O17").Value = 0.002553176504700935
$ws.Range("P17").Value = 0.002553176504700936
$ws.Range("Q17").Value = 0.2293206305973333
$ws.Range("R17").Value = 0.276200741327
$ws.Range("S17").Value = 0.00002511570813814356
$ws.Range("T17").Value = 0.00002511570813814356
